# "bug fix with import" - append the row that the importer produced for the
# newly-added person, reproducing the literal ("nothing1") uid value instead
# of a numeric id, and land the selection/scroll where the user ended up
# after entering the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A78").Value = "nothing1"
$ws.Range("B78").Value = "Чеченцев Герман Эльклассико"
$ws.Range("C78").Value = "ЭУ/с-21-1-о"

# Scroll the view down so row 65 is the first visible row, matching the
# position the user had scrolled to after adding the new row.
$win = $excel.ActiveWindow
$win.ScrollRow = 65
$win.ScrollColumn = 1

# Final selection lands a couple of rows below the new data, as in the
# source workbook.
$ws.Range("C80").Select()
